$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.603.25'
$ws.Range("E2").Value = '  +0.96%  '

$ws.Range("D3").Value = '3.878.89'
$ws.Range("E3").Value = '  +0.66%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'462.81"
$ws.Range("E5").Value = '  +9.18%  '

$ws.Range("D6").Value = "'148.87"
$ws.Range("E6").Value = '  +14.61%  '

$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = '  +3.30%  '

$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").Value = "'0.753"
$ws.Range("E9").Value = '  +4.14%  '

$ws.Range("D10").Value = "'0.158"
$ws.Range("E10").Value = '  +0.74%  '

$ws.Range("D11").Value = "'0.0000316"
$ws.Range("E11").Value = '  -6.01%  '

$ws.Range("D12").Value = "'44.20"
$ws.Range("E12").Value = '  +8.27%  '

$ws.Range("E13").Value = '  +1.17%  '

$ws.Range("D14").Value = '4.492.97'
$ws.Range("E14").Value = '  +0.80%  '

$ws.Range("D15").Value = "'14.86"
$ws.Range("E15").Value = '  -5.70%  '

$ws.Range("D16").Value = '3.831.55'
$ws.Range("E16").Value = '  -0.50%  '

$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("D18").Value = "'20.16"
$ws.Range("E18").Value = '  +1.59%  '

$ws.Range("E19").Value = '  +7.20%  '

$ws.Range("D20").Value = '67.700.88'
$ws.Range("E20").Value = '  +0.80%  '

$ws.Range("D21").Value = "'431.35"
$ws.Range("E21").Value = '  +4.24%  '

$ws.Range("D22").Value = "'14.93"
$ws.Range("E22").Value = '  -0.13%  '

$ws.Range("D23").Value = "'3.30"
$ws.Range("E23").Value = '  +8.50%  '

$ws.Range("D24").Value = "'88.02"
$ws.Range("E24").Value = '  +4.54%  '

$ws.Range("E25").Value = '  +9.65%  '

$ws.Range("D26").Value = "'10.37"
$ws.Range("E26").Value = '  +10.68%  '

$ws.Range("D27").Value = "'37.84"
$ws.Range("E27").Value = '  +0.96%  '

$ws.Range("D28").Value = "'10.19"
$ws.Range("E28").Value = '  +2.67%  '

$ws.Range("D29").Value = "'5.51"
$ws.Range("E29").Value = '  +2.99%  '

$ws.Range("D30").Value = "'751.13"
$ws.Range("E30").Value = '  +1.77%  '

$ws.Range("D31").Value = "'0.135"
$ws.Range("E31").Value = '  +10.83%  '

$ws.Range("D32").Value = "'13.77"
$ws.Range("E32").Value = '  +5.03%  '

$ws.Range("E33").Value = '  -0.54%  '

$ws.Range("D34").Value = "'43.66"
$ws.Range("E34").Value = '  +12.93%  '

$ws.Range("D35").Value = "'0.164"

$ws.Range("D36").Value = "'57.62"
$ws.Range("E36").Value = '  +3.83%  '

$ws.Range("E37").Value = '  +0.90%  '

$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("D39").Value = "'0.0480"
$ws.Range("E39").Value = '  +4.56%  '

$ws.Range("D40").Value = "'0.352"
$ws.Range("E40").Value = '  +11.10%  '

$ws.Range("D41").Value = "'2.95"
$ws.Range("E41").Value = '  +1.78%  '

$ws.Range("D42").Value = "'2.62"
$ws.Range("E42").Value = '  +13.86%  '

$ws.Range("E43").Value = '  -6.86%  '

$ws.Range("E44").Value = '  +5.42%  '

$ws.Range("E45").Value = '  +0.14%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = "'3.29"
$ws.Range("E46").Value = '  +6.22%  '

$ws.Range("B47").Value = 'LidoDAOToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D47").Value = "'3.45"
$ws.Range("E47").Value = '  +2.55%  '

$ws.Range("E48").Value = '  +7.67%  '

$ws.Range("E49").Value = '  +4.39%  '

$ws.Range("D50").Value = "'145.17"
$ws.Range("E50").Value = '  +3.18%  '

$ws.Range("D51").Value = "'2.91"
$ws.Range("E51").Value = '  +3.50%  '
